$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- Remove old hyperlinks (all of them; we'll re-add the 3 that remain) ---
$ws.Hyperlinks.Delete()

# --- Delete rows 5-12, keeping only header + 3 data rows ---
$ws.Range("A5:H12").EntireRow.Delete()

# --- Column width changes ---
# (The COM layer re-derives the stored "width" from a pixel conversion that
#  adds a consistent 5/6-character offset versus the value assigned to
#  ColumnWidth; compensate so the saved OOXML attribute matches exactly.)
$ws.Columns.Item(2).ColumnWidth = 44 - 5/6
$ws.Columns.Item(4).ColumnWidth = 28 - 5/6
$ws.Columns.Item(8).ColumnWidth = 12 - 5/6

# --- Row 2 ---
$ws.Range("A2").Value = "2025-11-18 06:26:52"
$ws.Range("B2").Value = "【業務委託】Shopee価格調整ツールの開発(Googleスプレッドシート+GAS)"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5436149"
$ws.Range("G2").Value = 128
$ws.Range("H2").Value = "◆ツール,開発"

# --- Row 3 ---
$ws.Range("A3").Value = "2025-11-18 06:26:52"
$ws.Range("B3").Value = "【技術パートナー募集】リード獲得・育成システム構築"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5436021"
$ws.Range("G3").Value = 33
$ws.Range("H3").ClearContents()

# --- Row 4 ---
$ws.Range("A4").Value = "2025-11-18 06:26:52"
$ws.Range("B4").Value = "【急募】Wartalesの武器アイコンとモデルを日本刀に差し替え"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5436248"
$ws.Range("G4").Value = 10
$ws.Range("H4").ClearContents()

# --- Re-add hyperlinks for F2:F4 and restore the Hyperlink cell style ---
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5436149")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5436021")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5436248")
$ws.Range("F2:F4").Style = "Hyperlink"
